$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (new weekly data point) ---
$ws.Range("D2").Value = 44993
$ws.Range("M2").Value = 14
$ws.Range("O2").Value = 200000
$ws.Range("P2").Value = 190000
$ws.Range("S2").Value = 190000

# --- Insert new row 4, a copy of the original row 2 values ---
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44672
$ws.Range("D4").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104005
$ws.Range("J4").Value = "Pera asiática"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 180000
$ws.Range("O4").Value = 180000
$ws.Range("P4").Value = 180000
$ws.Range("Q4").Value = "$/bins (430 kilos)"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 180000
$ws.Range("T4").Value = 1
